# Auto-generated Excel COM-interop edit script
# Applies the "Refactor application branding and enhance assessment features" changes
# to the High_Scoring_Candidates worksheet: adds Fundamental Knowledge Assessment
# columns (Q:T), normalizes row 2 (Gina Jones) name casing/content, and appends a
# new row 3 for candidate Ophelia Crane.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New header columns Q1:T1 - copy formatting (bold + border) from P1
# ---------------------------------------------------------------------------
$ws.Range("P1").Copy() | Out-Null
$ws.Range("Q1:T1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("Q1").Value = 'Fundamental Knowledge Score'
$ws.Range("R1").Value = 'Fundamental Recommendation'
$ws.Range("S1").Value = 'Fundamental Assessment Date'
$ws.Range("T1").Value = 'Fundamental Question Scores'

# ---------------------------------------------------------------------------
# 2) Row 2 (Gina Jones) - updated assessment text + new Fundamental columns
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = '2026-01-19 01:05 AM'

$ws.Range("C2").Value = 'Gina'

$ws.Range("D2").Value = 'Jones'

$v_F2 = @'
Practical experience in front-end interface development (8+ years at Adobe, Netflix, Airbnb).
Solid understanding of HTML5 (Implicit from extensive front-end and React development).
Solid understanding of CSS3 (Implicit from responsive UI design with Chakra UI and general front-end roles).
Solid understanding of JavaScript (ES6 basics) (Explicitly demonstrated by 'Certified Web Developer (JavaScript)' and extensive React/JS experience).
Proficiency with at least one modern front-end framework (e.g., React, Vue, Angular) OR a CSS framework (e.g., Bootstrap, Tailwind) combined with strong vanilla JavaScript (Extensive React experience and Chakra UI).
Basic familiarity with Git and common version control workflows (Implicit from working at major tech companies and 'Enforced code quality across the React codebase').
Ability to interpret design mockups (e.g., Figma, XD, or similar) and translate them into code ('Designed a responsive user interface using Chakra UI' and general UI development experience).
Ability to debug basic front-end issues using browser developer tools (Implicit from all front-end development roles and 'Enforced code quality').
'@
$ws.Range("F2").Value = $v_F2

$ws.Range("G2").Value = ''

$v_H2 = @'
Experience with React and its ecosystem (hooks, basic state management, routing) (Extensive React Developer experience, Redux, React Router).
Experience connecting to REST APIs and handling JSON data in complex flows (Used Apollo Client for GraphQL integration, implying strong API interaction skills; GraphQL often more complex than basic REST).
Exposure to performance optimization and accessibility best practices (Strong on performance optimization: 'optimizing application performance', 'reduced response time', 'decrease in bounce rate').
Basic front-end testing experience (e.g., Jest, React Testing Library, or similar) ('Created comprehensive test suites using Enzyme').
Experience working in a team environment (e.g., Agile/Scrum, code reviews) (Implicit from working at large tech companies; 'Enforced code quality across the React codebase with ESLint' implies code review involvement).
Familiarity with front-end build tools/bundlers (e.g., Webpack, Vite, npm/yarn) (Implicit from extensive modern JS framework development).
'@
$ws.Range("H2").Value = $v_H2

$ws.Range("I2").Value = 'Experience with TypeScript in front-end projects.'

$v_J2 = @'
Extensive experience (8+ years) at top-tier tech companies (Airbnb, Netflix, Adobe) in relevant front-end and full-stack roles.
Deep expertise in React and its ecosystem (Redux for state management, React Router, Enzyme for testing).
Proven track record of delivering measurable improvements in code quality (ESLint), test coverage, and application performance.
Experience with API integration (Apollo Client for GraphQL, implying general API understanding).
Strong understanding of UI development, including responsive design and component libraries (Chakra UI).
Certified Web Developer (JavaScript) demonstrating foundational knowledge.
'@
$ws.Range("J2").Value = $v_J2

$ws.Range("K2").Value = 'No explicit mention of TypeScript experience in front-end projects.'

$ws.Range("L2").Value = 'Low - Candidate has 8+ years of consistent experience at highly reputable tech companies (Adobe, Netflix, Airbnb) in roles directly relevant to front-end development. All core must-haves are met with significant depth, indicating a strong and reliable professional.'

$ws.Range("M2").Value = 'High - The candidate brings senior-level expertise, a proven track record of delivering measurable improvements (performance, test coverage, code quality), and extensive experience with modern front-end technologies like React. Their background suggests they can contribute significantly beyond the basic requirements and potentially take on leadership or mentorship roles.'

$ws.Range("N2").Value = 9

$ws.Range("O2").Value = 'Gina is an exceptionally strong candidate for the Front-End Developer role. She satisfies all must-have requirements, often demonstrating a level of experience and proficiency that significantly exceeds the basic ask. Her extensive background as a React Developer at Airbnb and Netflix, coupled with her Full-stack experience at Adobe, showcases deep expertise in modern front-end development, including React, state management (Redux), testing (Enzyme), and performance optimization. She also satisfies almost all nice-to-have requirements, with the only minor gap being explicit TypeScript experience. Her work at top-tier tech companies and measurable achievements (e.g., reducing code review time, increasing test coverage, optimizing performance) make her a high-reward, low-risk candidate who would likely excel in this role and contribute significantly from day one.'

$ws.Range("P2").Value = 'GINA JONES React Developerg.jones@email.com (123) 456-7890 San Francisco, CA LinkedIn Work Experience Airbnb React Developer 2020 - current|San Francisco, CA Enforced code quality across the React codebase with ESLint, resulting in a 21% reduction in code review time and improved overall codebase health. Created comprehensive test suites using Enzyme, which increased test coverage by 52%. Developed server-side applications using NestJS that reduced response time by 44%. Netﬂix Front-end Developer 2017 - 2020|Los Gatos, CA Integrated Redux for state management, optimizing application performance and reducing load time by 47%. Used Apollo Client for GraphQL integration, which resulted in a 53% reduction in API response time. Adobe Full-stack Developer 2015 - 2017|San Jose, CA Designed a responsive user interface using Chakra UI, which led to a 26% decrease in bounce rate. Optimized database queries to improve overall database performance by 43%. SkillsReact Router;ESLint;Chakra UI;Redux;Enzyme;Apollo Client;NestJS;Jekyll Education Bachelor of Science,Computer Science University of California 2011 - 2015|Berkeley, CA CertiﬁcationsCertiﬁed Web Developer (JavaScript)'

$ws.Range("Q2").Value = 17

$ws.Range("R2").Value = 'Not Recommended (Rejection)'

$ws.Range("S2").Value = '2026-01-19 01:09:51'

$v_T2 = @'
{"Q1: Write a JavaScript function to...": 100}
'@
$ws.Range("T2").Value = $v_T2

# ---------------------------------------------------------------------------
# 3) New row 3 (Ophelia Crane) - full new candidate record
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = '2026-01-19 01:06 AM'

$ws.Range("B3").Value = '1.pdf'

$ws.Range("C3").Value = 'Ophelia'

$ws.Range("D3").Value = 'Crane'

$ws.Range("E3").Value = 'o.crane@email.com'

$v_F3 = @'
Practical experience in front-end interface development (via work, internships, or significant personal/academic projects).
Solid understanding of HTML5 (semantic elements, forms, basic structure).
Solid understanding of CSS3 (layouts, positioning, responsive basics, fonts, colors).
Solid understanding of JavaScript (ES6 basics: variables, functions, arrays, events, simple DOM manipulation).
Proficiency with at least one modern front-end framework (e.g., React, Vue, Angular) OR a CSS framework (e.g., Bootstrap, Tailwind) combined with strong vanilla JavaScript.
Basic familiarity with Git and common version control workflows (clone, commit, push, pull).
Ability to interpret design mockups (e.g., Figma, XD, or similar) and translate them into code.
Ability to debug basic front-end issues using browser developer tools.
'@
$ws.Range("F3").Value = $v_F3

$ws.Range("G3").Value = ''

$v_H3 = @'
Experience with React and its ecosystem (hooks, basic state management, routing).
Exposure to performance optimization and accessibility best practices.
Basic front-end testing experience (e.g., Jest, React Testing Library, or similar).
Experience working in a team environment (e.g., Agile/Scrum, code reviews).
'@
$ws.Range("H3").Value = $v_H3

$v_I3 = @'
Experience with TypeScript in front-end projects.
Experience connecting to REST APIs and handling JSON data in complex flows.
Familiarity with front-end build tools/bundlers (e.g., Webpack, Vite, npm/yarn).
'@
$ws.Range("I3").Value = $v_I3

$v_J3 = @'
Strong foundational skills in HTML5, CSS3 (via Bootstrap and responsive design), React, and Git/GitHub.
Practical application of front-end skills demonstrated through significant academic projects and a non-dev role.
Experience with front-end testing (Jest) and working in a cross-functional team environment.
Proactive approach to problem-solving and improving user experience (e.g., reducing form errors, improving client satisfaction, enhancing accessibility).
Clear career objective aligned with a Front-End Developer role, indicating strong motivation.
'@
$ws.Range("J3").Value = $v_J3

$v_K3 = @'
Limited explicit professional full-time front-end development experience (primarily academic/project-based and an internship objective).
Lack of explicit experience with TypeScript in front-end projects.
No explicit mention of connecting to REST APIs and handling JSON data in complex flows.
No explicit familiarity with front-end build tools/bundlers (e.g., Webpack, Vite, npm/yarn).
'@
$ws.Range("K3").Value = $v_K3

$v_L3 = @'
Medium - The candidate's experience is primarily academic and project-based, with an objective for an internship. While skills are strong, there's a moderate risk regarding their ability to immediately handle complex production environments, advanced debugging, and full API integration without more direct professional experience. However, the job description is open to junior candidates, mitigating some of this risk.
'@
$ws.Range("L3").Value = $v_L3

$ws.Range("M3").Value = 'Medium - The candidate demonstrates strong initiative and a solid grasp of core front-end technologies. Their project work shows practical application, problem-solving abilities (e.g., reducing errors, improving satisfaction), and a results-oriented mindset. They possess key nice-to-have skills like React, Jest, and teamwork, indicating high growth potential for a junior role.'

$ws.Range("N3").Value = 8

$v_O3 = @'
The candidate satisfies all must-have requirements, demonstrating a solid foundation in HTML5, CSS3 (via Bootstrap and responsive design), JavaScript (implied by React), React, Git, and basic debugging/design interpretation. They also meet several nice-to-have criteria, including experience with React, front-end testing (Jest), accessibility concepts, and working in a team environment. While their experience is primarily academic and project-based, it is highly relevant and showcases practical application of skills. The candidate presents as a strong junior profile with significant growth potential, aligning well with the job description's openness to all experience levels. The moderate risk is balanced by their strong foundational skills and initiative.
'@
$ws.Range("O3").Value = $v_O3

$v_P3 = @'
OPHELIA CRANE Front End Developer Intern o.crane@email.com (123) 456-7890 Atlanta, GA LinkedIn EDUCATION Bachelor of Science Computer Science Georgia Institute of Technology 2021 - current Atlanta, GA SKILLS HTML5 React Git Bootstrap Jest GitHubCAREER OBJECTIVE Proactive front-end developer proﬁcient in HTML5 and React.js, seeking an internship position at Gretrix. I aim to apply my strong basics in front-end technologies and various languages to expedite the delivery of responsive and efﬁcient applications for every Gretrix client. W ORK EXPERIENCE Pet Sitter Critter Sitters Inc. 2023 - current Atlanta, GA Leveraged Git to manage version history of a new pet behavior tracking software, reducing average monthly incident reports by 34%. Implemented Bootstrap alerts for real-time updates on pet care status, improving communication and boosting client satisfaction scores by 32%. Created digital pet proﬁles, increasing information accessibility for clients and slashing support queries by 28%. Automated appointment reminders in advance via email, minimizing no- show rates by 14%. PROJECTS Library Assistant 2023 - 2024 Incorporated HTML5 form validation for the library's online registration system, lowering form submission errors by 42%. Developed a React-based notiﬁcation system for overdue books, decreasing overdue cases by 11%. Managed content updates on the library’s website, ensuring timely and accurate information dissemination. Worked with a cross-functional team to create an online study room reservation system, shrinking reservation conﬂicts by 26%. Campus Tour Guide 2022 Conducted 15+ campus tours weekly, providing informational experiences for prospective students, resulting in a 67% uplift in campus visit satisfaction scores. Built Jest test cases to verify the integration of multimedia elements in tour scripts, leading to a 21% surge in visitor engagement. Devised and maintained a GitHub repository for tour guide schedules and availability, curtailing the scheduling process time by 18 minutes. Revamped the site using responsive web design techniques to ensure the tour registration and information pages were accessible on all devices, growing mobile registrations by 28%.
'@
$ws.Range("P3").Value = $v_P3

$ws.Range("Q3").Value = 91

$ws.Range("R3").Value = 'Strong Yes (Excellent Candidate)'

$ws.Range("S3").Value = '2026-01-19 01:08:58'

$v_T3 = @'
{"Q1: Implement a responsive layout ...": 80}
'@
$ws.Range("T3").Value = $v_T3

Write-Output "High_Scoring_Candidates sheet updated: headers Q1:T1 added, row 2 refreshed, row 3 (Ophelia Crane) appended."